$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(97, 8).Value = 1040  # H97: 821.5714 -> 1040
$ws.Cells.Item(97, 10).Value = 1040  # J97: 821.5714 -> 1040
$ws.Cells.Item(97, 12).Value = 3120  # L97: 2464.7142 -> 3120
$ws.Cells.Item(97, 14).Value = -4112  # N97: -3456.7142 -> -4112

$ws.Cells.Item(113, 8).Value = 3598.3125  # H113: 2980.5908 -> 3598.3125
$ws.Cells.Item(113, 9).Value = 2529.2856  # I113: 2310.5 -> 2529.2856
$ws.Cells.Item(113, 10).Value = 4429.778  # J113: 3539 -> 4429.778
$ws.Cells.Item(113, 11).Value = 2529.2856  # K113: 2310.5 -> 2529.2856
$ws.Cells.Item(113, 12).Value = 4429.778  # L113: 3539 -> 4429.778
$ws.Cells.Item(113, 13).Value = 724.7143999999998  # M113: 943.5 -> 724.7143999999998
$ws.Cells.Item(113, 14).Value = -10937.778  # N113: -10047 -> -10937.778

$ws.Cells.Item(116, 8).Value = 7667.6665  # H116: 14003 -> 7667.6665
$ws.Cells.Item(116, 9).Value = 3000  # I116: 0 -> 3000
$ws.Cells.Item(116, 10).Value = 8601.200000000001  # J116: 14003 -> 8601.200000000001
$ws.Cells.Item(116, 11).Value = 3000  # K116: 0 -> 3000
$ws.Cells.Item(116, 12).Value = 8601.200000000001  # L116: 14003 -> 8601.200000000001
$ws.Cells.Item(116, 13).Value = 442  # M116: None -> 442
$ws.Cells.Item(116, 14).Value = -15485.2  # N116: -20887 -> -15485.2

$ws.Cells.Item(129, 8).Value = 200840.06  # H129: 223065.28 -> 200840.06
$ws.Cells.Item(129, 10).Value = 200840.06  # J129: 223065.28 -> 200840.06
$ws.Cells.Item(129, 12).Value = 602520.1799999999  # L129: 669195.84 -> 602520.1799999999
$ws.Cells.Item(129, 14).Value = -612520.1799999999  # N129: -679195.84 -> -612520.1799999999

$ws.Cells.Item(132, 8).Value = 7595.4546  # H132: 11171.429 -> 7595.4546
$ws.Cells.Item(132, 9).Value = 8061.1113  # I132: 13440 -> 8061.1113
$ws.Cells.Item(132, 11).Value = 24183.3339  # K132: 40320 -> 24183.3339
$ws.Cells.Item(132, 13).Value = -21653.3339  # M132: -37790 -> -21653.3339

$ws.Cells.Item(137, 8).Value = 78838.30499999999  # H137: 2066.5833 -> 78838.30499999999
$ws.Cells.Item(137, 9).Value = 1899.8  # I137: 1733.1666 -> 1899.8
$ws.Cells.Item(137, 10).Value = 126924.875  # J137: 2400 -> 126924.875
$ws.Cells.Item(137, 11).Value = 5699.4  # K137: 5199.4998 -> 5699.4
$ws.Cells.Item(137, 12).Value = 380774.625  # L137: 7200 -> 380774.625
$ws.Cells.Item(137, 13).Value = -3149.4  # M137: -2649.4998 -> -3149.4
$ws.Cells.Item(137, 14).Value = -385874.625  # N137: -12300 -> -385874.625

$ws.Cells.Item(138, 8).Value = 1800.5605  # H138: 1806.0156 -> 1800.5605
$ws.Cells.Item(138, 10).Value = 2383.647  # J138: 2431 -> 2383.647
$ws.Cells.Item(138, 12).Value = 7150.941  # L138: 7293 -> 7150.941
$ws.Cells.Item(138, 14).Value = -17430.941  # N138: -17573 -> -17430.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 0  # H6: 6999 -> 0
$ws.Cells.Item(6, 9).Value = 0  # I6: 6999 -> 0
$ws.Cells.Item(6, 11).Value = 0  # K6: 6999 -> 0
$ws.Cells.Item(6, 13).ClearContents()  # M6: remove (was -6826)

$ws.Cells.Item(19, 8).Value = 769.3333  # H19: 0 -> 769.3333
$ws.Cells.Item(19, 9).Value = 654  # I19: 0 -> 654
$ws.Cells.Item(19, 10).Value = 1000  # J19: 0 -> 1000
$ws.Cells.Item(19, 11).Value = 654  # K19: 0 -> 654
$ws.Cells.Item(19, 12).Value = 1000  # L19: 0 -> 1000
$ws.Cells.Item(19, 13).Value = -425  # M19: None -> -425
$ws.Cells.Item(19, 14).Value = -1458  # N19: None -> -1458

$ws.Cells.Item(26, 8).Value = 3081.4  # H26: 4500 -> 3081.4
$ws.Cells.Item(26, 9).Value = 3081.4  # I26: 4500 -> 3081.4
$ws.Cells.Item(26, 11).Value = 3081.4  # K26: 4500 -> 3081.4
$ws.Cells.Item(26, 13).Value = -2751.4  # M26: -4170 -> -2751.4

$ws.Cells.Item(39, 8).Value = 1000  # H39: 1016 -> 1000
$ws.Cells.Item(39, 9).Value = 1000  # I39: 1016 -> 1000
$ws.Cells.Item(39, 11).Value = 1000  # K39: 1016 -> 1000
$ws.Cells.Item(39, 13).Value = -480  # M39: -496 -> -480

$ws.Cells.Item(74, 8).Value = 66670108  # H74: 62503276 -> 66670108
$ws.Cells.Item(74, 9).Value = 71432150  # I74: 66670060 -> 71432150
$ws.Cells.Item(74, 11).Value = 71432150  # K74: 66670060 -> 71432150
$ws.Cells.Item(74, 13).Value = -71431276  # M74: -66669186 -> -71431276

$ws.Cells.Item(77, 8).Value = 66670108  # H77: 62503276 -> 66670108
$ws.Cells.Item(77, 9).Value = 71432150  # I77: 66670060 -> 71432150
$ws.Cells.Item(77, 11).Value = 357160750  # K77: 333350300 -> 357160750
$ws.Cells.Item(77, 13).Value = -357156382  # M77: -333345932 -> -357156382

$ws.Cells.Item(122, 8).Value = 1923.8636  # H122: 2085.7896 -> 1923.8636
$ws.Cells.Item(122, 9).Value = 2116.4375  # I122: 2385.8462 -> 2116.4375
$ws.Cells.Item(122, 10).Value = 1410.3334  # J122: 1435.6666 -> 1410.3334
$ws.Cells.Item(122, 11).Value = 6349.3125  # K122: 7157.5386 -> 6349.3125
$ws.Cells.Item(122, 12).Value = 4231.0002  # L122: 4306.9998 -> 4231.0002
$ws.Cells.Item(122, 13).Value = -3899.3125  # M122: -4707.5386 -> -3899.3125
$ws.Cells.Item(122, 14).Value = -9131.0002  # N122: -9206.9998 -> -9131.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 5002575  # H7: 3336766.2 -> 5002575
$ws.Cells.Item(7, 10).Value = 5005000  # J7: 10000 -> 5005000
$ws.Cells.Item(7, 12).Value = 5005000  # L7: 10000 -> 5005000
$ws.Cells.Item(7, 14).Value = -5005226  # N7: -10226 -> -5005226

$ws.Cells.Item(20, 8).Value = 3566.5833  # H20: 3800.0908 -> 3566.5833
$ws.Cells.Item(20, 9).Value = 4475.125  # I20: 4971.857 -> 4475.125
$ws.Cells.Item(20, 11).Value = 4475.125  # K20: 4971.857 -> 4475.125
$ws.Cells.Item(20, 13).Value = -4228.125  # M20: -4724.857 -> -4228.125

$ws.Cells.Item(94, 8).Value = 1063.8654  # H94: 1093.66 -> 1063.8654
$ws.Cells.Item(94, 9).Value = 916.2143  # I94: 946.075 -> 916.2143
$ws.Cells.Item(94, 11).Value = 916.2143  # K94: 946.075 -> 916.2143
$ws.Cells.Item(94, 13).Value = -465.2143  # M94: -495.075 -> -465.2143

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 294.5  # H12: 320 -> 294.5
$ws.Cells.Item(12, 9).Value = 289  # I12: 230 -> 289
$ws.Cells.Item(12, 10).Value = 300  # J12: 500 -> 300
$ws.Cells.Item(12, 11).Value = 289  # K12: 230 -> 289
$ws.Cells.Item(12, 12).Value = 300  # L12: 500 -> 300
$ws.Cells.Item(12, 13).Value = -119  # M12: -60 -> -119
$ws.Cells.Item(12, 14).Value = -640  # N12: -840 -> -640

$ws.Cells.Item(31, 8).Value = 22475.062  # H31: 3800 -> 22475.062
$ws.Cells.Item(31, 9).Value = 37455.668  # I31: 4030 -> 37455.668
$ws.Cells.Item(31, 10).Value = 3214.2856  # J31: 3416.6667 -> 3214.2856
$ws.Cells.Item(31, 11).Value = 37455.668  # K31: 4030 -> 37455.668
$ws.Cells.Item(31, 12).Value = 3214.2856  # L31: 3416.6667 -> 3214.2856
$ws.Cells.Item(31, 13).Value = -37160.668  # M31: -3735 -> -37160.668
$ws.Cells.Item(31, 14).Value = -3804.2856  # N31: -4006.6667 -> -3804.2856

$ws.Cells.Item(32, 8).Value = 2670  # H32: 2833 -> 2670
$ws.Cells.Item(32, 9).Value = 2755  # I32: 2833 -> 2755
$ws.Cells.Item(32, 10).Value = 2500  # J32: 0 -> 2500
$ws.Cells.Item(32, 11).Value = 2755  # K32: 2833 -> 2755
$ws.Cells.Item(32, 12).Value = 2500  # L32: 0 -> 2500
$ws.Cells.Item(32, 13).Value = -2439  # M32: -2517 -> -2439
$ws.Cells.Item(32, 14).Value = -3132  # N32: None -> -3132

$ws.Cells.Item(34, 8).Value = 22475.062  # H34: 3800 -> 22475.062
$ws.Cells.Item(34, 9).Value = 37455.668  # I34: 4030 -> 37455.668
$ws.Cells.Item(34, 10).Value = 3214.2856  # J34: 3416.6667 -> 3214.2856
$ws.Cells.Item(34, 11).Value = 37455.668  # K34: 4030 -> 37455.668
$ws.Cells.Item(34, 12).Value = 3214.2856  # L34: 3416.6667 -> 3214.2856
$ws.Cells.Item(34, 13).Value = -37253.668  # M34: -3828 -> -37253.668
$ws.Cells.Item(34, 14).Value = -3618.2856  # N34: -3820.6667 -> -3618.2856

$ws.Cells.Item(35, 8).Value = 2420.8333  # H35: 3500 -> 2420.8333
$ws.Cells.Item(35, 9).Value = 2420.8333  # I35: 3500 -> 2420.8333
$ws.Cells.Item(35, 11).Value = 2420.8333  # K35: 3500 -> 2420.8333
$ws.Cells.Item(35, 13).Value = -2126.8333  # M35: -3206 -> -2126.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 114419.375  # H131: 106046.21 -> 114419.375
$ws.Cells.Item(131, 9).Value = 900  # I131: 965 -> 900
$ws.Cells.Item(131, 10).Value = 115724.195  # J131: 108306.02 -> 115724.195
$ws.Cells.Item(131, 11).Value = 2700  # K131: 2895 -> 2700
$ws.Cells.Item(131, 12).Value = 347172.585  # L131: 324918.06 -> 347172.585
$ws.Cells.Item(131, 13).Value = 2340  # M131: 2145 -> 2340
$ws.Cells.Item(131, 14).Value = -357252.585  # N131: -334998.06 -> -357252.585

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 7000000  # H12: 6933333.5 -> 7000000
$ws.Cells.Item(12, 10).Value = 0  # J12: 6000000 -> 0
$ws.Cells.Item(12, 12).Value = 0  # L12: 6000000 -> 0
$ws.Cells.Item(12, 14).ClearContents()  # N12: remove (was -6000280)

$ws.Cells.Item(70, 8).Value = 5221217  # H70: 11001.643 -> 5221217
$ws.Cells.Item(70, 9).Value = 29925  # I70: 15898.375 -> 29925
$ws.Cells.Item(70, 10).Value = 7816863  # J70: 4472.6665 -> 7816863
$ws.Cells.Item(70, 11).Value = 29925  # K70: 15898.375 -> 29925
$ws.Cells.Item(70, 12).Value = 7816863  # L70: 4472.6665 -> 7816863
$ws.Cells.Item(70, 13).Value = -29655  # M70: -15628.375 -> -29655
$ws.Cells.Item(70, 14).Value = -7817403  # N70: -5012.6665 -> -7817403

$ws.Cells.Item(73, 8).Value = 5221217  # H73: 11001.643 -> 5221217
$ws.Cells.Item(73, 9).Value = 29925  # I73: 15898.375 -> 29925
$ws.Cells.Item(73, 10).Value = 7816863  # J73: 4472.6665 -> 7816863
$ws.Cells.Item(73, 11).Value = 29925  # K73: 15898.375 -> 29925
$ws.Cells.Item(73, 12).Value = 7816863  # L73: 4472.6665 -> 7816863
$ws.Cells.Item(73, 13).Value = -28989  # M73: -14962.375 -> -28989
$ws.Cells.Item(73, 14).Value = -7818735  # N73: -6344.6665 -> -7818735

$ws.Cells.Item(126, 8).Value = 4058.8235  # H126: 3712.1282 -> 4058.8235
$ws.Cells.Item(126, 9).Value = 3147.6191  # I126: 2802.8076 -> 3147.6191
$ws.Cells.Item(126, 11).Value = 9442.8573  # K126: 8408.4228 -> 9442.8573
$ws.Cells.Item(126, 13).Value = -6972.8573  # M126: -5938.4228 -> -6972.8573

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 3656.5  # H32: 3256.5 -> 3656.5
$ws.Cells.Item(32, 9).Value = 3656.5  # I32: 3256.5 -> 3656.5
$ws.Cells.Item(32, 11).Value = 3656.5  # K32: 3256.5 -> 3656.5
$ws.Cells.Item(32, 13).Value = -3339.5  # M32: -2939.5 -> -3339.5

$ws.Cells.Item(122, 8).Value = 1310783.5  # H122: 1229078.2 -> 1310783.5
$ws.Cells.Item(122, 10).Value = 5050.6665  # J122: 4829.143 -> 5050.6665
$ws.Cells.Item(122, 12).Value = 15151.9995  # L122: 14487.429 -> 15151.9995
$ws.Cells.Item(122, 14).Value = -20051.9995  # N122: -19387.429 -> -20051.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 395  # H23: 1000 -> 395
$ws.Cells.Item(23, 9).Value = 395  # I23: 1000 -> 395
$ws.Cells.Item(23, 11).Value = 395  # K23: 1000 -> 395
$ws.Cells.Item(23, 13).Value = -166  # M23: -771 -> -166

$ws.Cells.Item(29, 8).Value = 200  # H29: 0 -> 200
$ws.Cells.Item(29, 9).Value = 200  # I29: 0 -> 200
$ws.Cells.Item(29, 11).Value = 200  # K29: 0 -> 200
$ws.Cells.Item(29, 13).Value = 90  # M29: None -> 90

$ws.Cells.Item(126, 8).Value = 1232.4231  # H126: 1356 -> 1232.4231
$ws.Cells.Item(126, 9).Value = 1339.5625  # I126: 1666.2727 -> 1339.5625
$ws.Cells.Item(126, 10).Value = 1061  # J126: 1045.7273 -> 1061
$ws.Cells.Item(126, 11).Value = 4018.6875  # K126: 4998.8181 -> 4018.6875
$ws.Cells.Item(126, 12).Value = 3183  # L126: 3137.1819 -> 3183
$ws.Cells.Item(126, 13).Value = -1548.6875  # M126: -2528.8181 -> -1548.6875
$ws.Cells.Item(126, 14).Value = -8123  # N126: -8077.1819 -> -8123
